# This script applies an update to the NATMI TPM output table (Hras-Insr.xlsx).
# The underlying TPM input data changed, which altered several derived
# ligand/receptor expression and specificity columns (G,H,I,J,M,N,O,P,Q,R,S,T).
# Only the cells whose computed values actually changed are updated here,
# using the exact new values reproduced from the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 7.533107333333334
$ws.Cells.Item(2, 8).Value = 22.599322
$ws.Cells.Item(2, 9).Value = 0.4772251808959424
$ws.Cells.Item(2, 10).Value = 0.4772251808959424
$ws.Cells.Item(2, 13).Value = 8.533046666666666
$ws.Cells.Item(2, 14).Value = 25.59914
$ws.Cells.Item(2, 15).Value = 0.2932132236642383
$ws.Cells.Item(2, 16).Value = 0.2932132236642383
$ws.Cells.Item(2, 17).Value = 64.28035642034222
$ws.Cells.Item(2, 18).Value = 578.52320778308
$ws.Cells.Item(2, 19).Value = 0.1399287337042485
$ws.Cells.Item(2, 20).Value = 0.1399287337042485

$ws.Cells.Item(3, 7).Value = 7.533107333333334
$ws.Cells.Item(3, 8).Value = 22.599322
$ws.Cells.Item(3, 9).Value = 0.4772251808959424
$ws.Cells.Item(3, 10).Value = 0.4772251808959424
$ws.Cells.Item(3, 15).Value = 0.3119288965200195
$ws.Cells.Item(3, 16).Value = 0.3119288965200194
$ws.Cells.Item(3, 17).Value = 68.38334368258712
$ws.Cells.Item(3, 18).Value = 615.4500931432841
$ws.Cells.Item(3, 19).Value = 0.148860324068438
$ws.Cells.Item(3, 20).Value = 0.148860324068438

$ws.Cells.Item(4, 7).Value = 7.533107333333334
$ws.Cells.Item(4, 8).Value = 22.599322
$ws.Cells.Item(4, 9).Value = 0.4772251808959424
$ws.Cells.Item(4, 10).Value = 0.4772251808959424
$ws.Cells.Item(4, 15).Value = 0.3948578798157423
$ws.Cells.Item(4, 16).Value = 0.3948578798157423
$ws.Cells.Item(4, 17).Value = 86.56364447942266
$ws.Cells.Item(4, 18).Value = 779.072800314804
$ws.Cells.Item(4, 19).Value = 0.1884361231232559
$ws.Cells.Item(4, 20).Value = 0.1884361231232559

$ws.Cells.Item(5, 7).Value = 5.009378000000001
$ws.Cells.Item(5, 9).Value = 0.3173459790819593
$ws.Cells.Item(5, 10).Value = 0.3173459790819593
$ws.Cells.Item(5, 13).Value = 8.533046666666666
$ws.Cells.Item(5, 14).Value = 25.59914
$ws.Cells.Item(5, 15).Value = 0.2932132236642383
$ws.Cells.Item(5, 16).Value = 0.2932132236642383
$ws.Cells.Item(5, 17).Value = 42.74525624497333
$ws.Cells.Item(5, 18).Value = 384.70730620476
$ws.Cells.Item(5, 19).Value = 0.09305003754350523
$ws.Cells.Item(5, 20).Value = 0.09305003754350522

$ws.Cells.Item(6, 7).Value = 5.009378000000001
$ws.Cells.Item(6, 9).Value = 0.3173459790819593
$ws.Cells.Item(6, 10).Value = 0.3173459790819593
$ws.Cells.Item(6, 15).Value = 0.3119288965200195
$ws.Cells.Item(6, 16).Value = 0.3119288965200194
$ws.Cells.Item(6, 17).Value = 45.47366740603868
$ws.Cells.Item(6, 18).Value = 409.2630066543481
$ws.Cells.Item(6, 19).Value = 0.09898938107010076
$ws.Cells.Item(6, 20).Value = 0.09898938107010072

$ws.Cells.Item(7, 7).Value = 5.009378000000001
$ws.Cells.Item(7, 9).Value = 0.3173459790819593
$ws.Cells.Item(7, 10).Value = 0.3173459790819593
$ws.Cells.Item(7, 15).Value = 0.3948578798157423
$ws.Cells.Item(7, 16).Value = 0.3948578798157423
$ws.Cells.Item(7, 17).Value = 57.563233479532
$ws.Cells.Item(7, 18).Value = 518.069101315788
$ws.Cells.Item(7, 19).Value = 0.1253065604683534
$ws.Cells.Item(7, 20).Value = 0.1253065604683533

$ws.Cells.Item(8, 9).Value = 0.2054288400220983
$ws.Cells.Item(8, 10).Value = 0.2054288400220983
$ws.Cells.Item(8, 13).Value = 8.533046666666666
$ws.Cells.Item(8, 14).Value = 25.59914
$ws.Cells.Item(8, 15).Value = 0.2932132236642383
$ws.Cells.Item(8, 16).Value = 0.2932132236642383
$ws.Cells.Item(8, 17).Value = 27.67045743656444
$ws.Cells.Item(8, 18).Value = 249.03411692908
$ws.Cells.Item(8, 19).Value = 0.06023445241648454
$ws.Cells.Item(8, 20).Value = 0.06023445241648454

$ws.Cells.Item(9, 9).Value = 0.2054288400220983
$ws.Cells.Item(9, 10).Value = 0.2054288400220983
$ws.Cells.Item(9, 15).Value = 0.3119288965200195
$ws.Cells.Item(9, 16).Value = 0.3119288965200194
$ws.Cells.Item(9, 19).Value = 0.06407919138148074
$ws.Cells.Item(9, 20).Value = 0.06407919138148073

$ws.Cells.Item(10, 9).Value = 0.2054288400220983
$ws.Cells.Item(10, 10).Value = 0.2054288400220983
$ws.Cells.Item(10, 15).Value = 0.3948578798157423
$ws.Cells.Item(10, 16).Value = 0.3948578798157423
$ws.Cells.Item(10, 19).Value = 0.08111519622413303
$ws.Cells.Item(10, 20).Value = 0.08111519622413303
